$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newUrl  = "https://www.london-fire.gov.uk/incidents/2021/may/industrial-unit-fire-hampton/"
$newText = "Industrial unit fire - Hampton | London Fire Brigade (london-fire.gov.uk)"

# Add the new source-link row for the Hampton industrial unit fire, mirroring
# the other link cells already on the sheet (D5/E5/F3/E4/E2): the hyperlink's
# cached display text is the raw URL, while the cell itself shows the
# friendly article title.
$ws.Hyperlinks.Add($ws.Range("F5"), $newUrl, "", "", $newUrl) | Out-Null
$ws.Range("F5").Value = $newText

# Hyperlinks.Add() stamps its own Hyperlink style instance onto the cell;
# reassign the named style so F5 matches the workbook's existing hyperlink
# cells exactly.
$ws.Range("F5").Style = "Hyperlink"

# Leave the cursor on the newly-added cell, like the editing user would.
$ws.Range("F5").Select() | Out-Null
